$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2 through 297 all currently hold the date
# serial value 45188; update them to 45189.
$ws.Range("C2:C297").Value = 45189
